$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1344.6364
$ws.Range("I2").Value = 698.4
$ws.Range("J2").Value = 1883.1666
$ws.Range("K2").Value = 698.4
$ws.Range("L2").Value = 1883.1666
$ws.Range("M2").Value = -585.4
$ws.Range("N2").Value = -2109.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6122.8887
$ws.Range("I19").Value = 5123.5
$ws.Range("J19").Value = 6408.4287
$ws.Range("K19").Value = 5123.5
$ws.Range("L19").Value = 6408.4287
$ws.Range("M19").Value = -4948.5
$ws.Range("N19").Value = -6758.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5141.8184
$ws.Range("J40").Value = 5991.25
$ws.Range("L40").Value = 5991.25
$ws.Range("N40").Value = -6341.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11887.223
$ws.Range("I43").Value = 16750
$ws.Range("K43").Value = 16750
$ws.Range("M43").Value = -16681

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6491.3125
$ws.Range("I70").Value = 992.3333
$ws.Range("J70").Value = 7760.3076
$ws.Range("K70").Value = 2976.9999
$ws.Range("L70").Value = 23280.9228
$ws.Range("M70").Value = -2706.9999
$ws.Range("N70").Value = -23820.9228

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6491.3125
$ws.Range("I73").Value = 992.3333
$ws.Range("J73").Value = 7760.3076
$ws.Range("K73").Value = 2976.9999
$ws.Range("L73").Value = 23280.9228
$ws.Range("M73").Value = -2040.9999
$ws.Range("N73").Value = -25152.9228

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4253.4375
$ws.Range("I76").Value = 4075.3572
$ws.Range("K76").Value = 4075.3572
$ws.Range("M76").Value = -3760.3572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4253.4375
$ws.Range("I79").Value = 4075.3572
$ws.Range("K79").Value = 4075.3572
$ws.Range("M79").Value = -2983.3572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 6290.16
$ws.Range("I80").Value = 483.2
$ws.Range("K80").Value = 1449.6
$ws.Range("M80").Value = -451.5999999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 6290.16
$ws.Range("I83").Value = 483.2
$ws.Range("K83").Value = 4348.8
$ws.Range("M83").Value = 643.1999999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6091.5
$ws.Range("I86").Value = 2869.4285
$ws.Range("J86").Value = 10602.4
$ws.Range("K86").Value = 2869.4285
$ws.Range("L86").Value = 10602.4
$ws.Range("M86").Value = -1746.4285
$ws.Range("N86").Value = -12848.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 12249
$ws.Range("J88").Value = 13427.429
$ws.Range("L88").Value = 13427.429
$ws.Range("N88").Value = -14239.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6091.5
$ws.Range("I89").Value = 2869.4285
$ws.Range("J89").Value = 10602.4
$ws.Range("K89").Value = 14347.1425
$ws.Range("L89").Value = 53012
$ws.Range("M89").Value = -8731.1425
$ws.Range("N89").Value = -64244

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 12249
$ws.Range("J91").Value = 13427.429
$ws.Range("L91").Value = 13427.429
$ws.Range("N91").Value = -16235.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4182.1
$ws.Range("I113").Value = 3845
$ws.Range("J113").Value = 4687.75
$ws.Range("K113").Value = 3845
$ws.Range("L113").Value = 4687.75
$ws.Range("M113").Value = -591
$ws.Range("N113").Value = -11195.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2816134.5
$ws.Range("I132").Value = 3263941.2
$ws.Range("J132").Value = 1348.4286
$ws.Range("K132").Value = 9791823.600000001
$ws.Range("L132").Value = 4045.2858
$ws.Range("M132").Value = -9789293.600000001
$ws.Range("N132").Value = -9105.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 28207
$ws.Range("I135").Value = 992.8570999999999
$ws.Range("K135").Value = 8935.713899999999
$ws.Range("M135").Value = -6400.713899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20096.219
$ws.Range("I32").Value = 21484.646
$ws.Range("K32").Value = 21484.646
$ws.Range("M32").Value = -21197.646

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 222348.08
$ws.Range("I74").Value = 273657.2
$ws.Range("J74").Value = 34214.668
$ws.Range("K74").Value = 273657.2
$ws.Range("L74").Value = 34214.668
$ws.Range("M74").Value = -272783.2
$ws.Range("N74").Value = -35962.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 222348.08
$ws.Range("I77").Value = 273657.2
$ws.Range("J77").Value = 34214.668
$ws.Range("K77").Value = 1368286
$ws.Range("L77").Value = 171073.34
$ws.Range("M77").Value = -1363918
$ws.Range("N77").Value = -179809.34

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1616.1666
$ws.Range("I97").Value = 1201.3914
$ws.Range("K97").Value = 1201.3914
$ws.Range("M97").Value = -705.3914

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 887.53845
$ws.Range("I94").Value = 909.2941
$ws.Range("K94").Value = 909.2941
$ws.Range("M94").Value = -458.2941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 47811.09
$ws.Range("I107").Value = 71958.42999999999
$ws.Range("J107").Value = 5553.25
$ws.Range("K107").Value = 71958.42999999999
$ws.Range("L107").Value = 5553.25
$ws.Range("M107").Value = -70038.42999999999
$ws.Range("N107").Value = -9393.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 23266.2
$ws.Range("J55").Value = 34110.332
$ws.Range("L55").Value = 34110.332
$ws.Range("N55").Value = -34740.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2418.2856
$ws.Range("I94").Value = 2600
$ws.Range("J94").Value = 2388
$ws.Range("K94").Value = 2600
$ws.Range("L94").Value = 2388
$ws.Range("M94").Value = -2149
$ws.Range("N94").Value = -3290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 42717.668
$ws.Range("I132").Value = 42717.668
$ws.Range("K132").Value = 128153.004
$ws.Range("M132").Value = -125623.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4999.975
$ws.Range("J81").Value = 4999.975
$ws.Range("L81").Value = 14999.925
$ws.Range("N81").Value = -17245.925

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 4999.975
$ws.Range("J84").Value = 4999.975
$ws.Range("L84").Value = 44999.775
$ws.Range("N84").Value = -56231.775

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 1992.625
$ws.Range("I123").Value = 988.6
$ws.Range("J123").Value = 3666
$ws.Range("K123").Value = 2965.8
$ws.Range("L123").Value = 10998
$ws.Range("M123").Value = -515.8000000000002
$ws.Range("N123").Value = -15898

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 1905
$ws.Range("I126").Value = 1905
$ws.Range("K126").Value = 5715
$ws.Range("M126").Value = -775

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 1765
$ws.Range("J127").Value = 2500
$ws.Range("L127").Value = 7500
$ws.Range("N127").Value = -17420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18440.625
$ws.Range("I15").Value = 16000
$ws.Range("J15").Value = 18789.285
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 18789.285
$ws.Range("M15").Value = -15712
$ws.Range("N15").Value = -19365.285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7457.0938
$ws.Range("J70").Value = 7768.1816
$ws.Range("L70").Value = 7768.1816
$ws.Range("N70").Value = -8308.1816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7457.0938
$ws.Range("J73").Value = 7768.1816
$ws.Range("L73").Value = 7768.1816
$ws.Range("N73").Value = -9640.1816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 18440.625
$ws.Range("I81").Value = 16000
$ws.Range("J81").Value = 18789.285
$ws.Range("K81").Value = 16000
$ws.Range("L81").Value = 18789.285
$ws.Range("M81").Value = -15002
$ws.Range("N81").Value = -20785.285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 18440.625
$ws.Range("I84").Value = 16000
$ws.Range("J84").Value = 18789.285
$ws.Range("K84").Value = 48000
$ws.Range("L84").Value = 56367.855
$ws.Range("M84").Value = -43008
$ws.Range("N84").Value = -66351.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2762
$ws.Range("J97").Value = 2961.1538
$ws.Range("L97").Value = 2961.1538
$ws.Range("N97").Value = -3953.1538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3471.9048
$ws.Range("I122").Value = 3516.2307
$ws.Range("J122").Value = 3399.875
$ws.Range("K122").Value = 10548.6921
$ws.Range("L122").Value = 10199.625
$ws.Range("M122").Value = -8098.6921
$ws.Range("N122").Value = -15099.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3638.1538
$ws.Range("I126").Value = 2216.5
$ws.Range("K126").Value = 6649.5
$ws.Range("M126").Value = -4179.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2755.7144
$ws.Range("I7").Value = 2715.0833
$ws.Range("K7").Value = 2715.0833
$ws.Range("M7").Value = -2603.0833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3869.1428
$ws.Range("I122").Value = 2476.2856
$ws.Range("K122").Value = 7428.8568
$ws.Range("M122").Value = -4978.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2755.7144
$ws.Range("I126").Value = 2715.0833
$ws.Range("K126").Value = 8145.249899999999
$ws.Range("M126").Value = -5675.249899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 19149.615
$ws.Range("J81").Value = 6000
$ws.Range("L81").Value = 12000
$ws.Range("N81").Value = -14122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 19149.615
$ws.Range("J84").Value = 6000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -70608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1219.8
$ws.Range("I96").Value = 1219.8
$ws.Range("K96").Value = 1219.8
$ws.Range("M96").Value = 153.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 703.5714
$ws.Range("I107").Value = 629
$ws.Range("K107").Value = 1887
$ws.Range("M107").Value = 33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 154510.45
$ws.Range("I126").Value = 2344.9524
$ws.Range("J126").Value = 420800.1
$ws.Range("K126").Value = 7034.8572
$ws.Range("L126").Value = 1262400.3
$ws.Range("M126").Value = -4564.8572
$ws.Range("N126").Value = -1267340.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10811.121
$ws.Range("I136").Value = 11158.352
$ws.Range("K136").Value = 33475.056
$ws.Range("M136").Value = -30925.056
